# Automatic data refresh for resum_periode_meteocat.xlsx
# Commit: "Update automàtic: dades i banners [2026-02-20 11:15]"
#
# The workbook stores every measurement as plain text (inlineStr in the
# original file), even values that happen to look numeric (e.g. "194",
# "15.0"). To keep those as text through COM automation -- instead of
# having Excel silently coerce them into numbers -- a leading apostrophe
# is used, exactly as a user typing into Excel would do to force text.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Dades_Període")
$ws2 = $wb.Worksheets.Item("Estudi_Capçaleres")

# ---------------------------------------------------------------------
# Sheet "Dades_Període" - row 2 (current 10:00-10:30 period becomes the
# newly observed 10:30-11:00 period)
# ---------------------------------------------------------------------

# Plain text values (not numeric-looking, no quoting needed)
$ws1.Range("E2").Value = "10:30 - 11:00"
$ws1.Range("H2").Value = "2026-02-20 11:15:27"
$ws1.Range("I2").Value = "10:30"
$ws1.Range("J2").Value = "https://www.meteo.cat/observacions/xema/dades?codi=XJ&dia=2026-02-20T10:30Z"
$ws1.Range("Q2").Value = "10:30 - 11:00"
$ws1.Range("X2").Value = "10:30 - 11:00"
$ws1.Range("AI2").Value = "10:30 - 11:00"

# Numeric-looking values that must remain stored as text
$ws1.Range("M2").Value = "'194"
$ws1.Range("N2").Value = "'42"

$ws1.Range("R2").Value = "'649"
$ws1.Range("S2").Value = "'15.0"
$ws1.Range("T2").Value = "'14.4"
$ws1.Range("U2").Value = "'15.6"
$ws1.Range("V2").Value = "'4.3"
$ws1.Range("W2").Value = "'11.9"

$ws1.Range("Y2").Value = "'15.0"
$ws1.Range("Z2").Value = "'15.6"
$ws1.Range("AA2").Value = "'14.4"
$ws1.Range("AB2").Value = "'42"
$ws1.Range("AD2").Value = "'4.3"
$ws1.Range("AE2").Value = "'194"
$ws1.Range("AF2").Value = "'11.9"
$ws1.Range("AH2").Value = "'649"

$ws1.Range("AJ2").Value = "'15.0"
$ws1.Range("AK2").Value = "'15.6"
$ws1.Range("AL2").Value = "'14.4"
$ws1.Range("AM2").Value = "'42"
$ws1.Range("AO2").Value = "'4.3"
$ws1.Range("AP2").Value = "'194"
$ws1.Range("AQ2").Value = "'11.9"
$ws1.Range("AS2").Value = "'649"

# ---------------------------------------------------------------------
# Sheet "Dades_Període" - DATA_EXTRACCIO (column H) refreshed for the
# other rows too (same extraction pass, new timestamp)
# ---------------------------------------------------------------------
$ws1.Range("H3").Value = "2026-02-20 11:15:28"
$ws1.Range("H4").Value = "2026-02-20 11:15:28"
$ws1.Range("H5").Value = "2026-02-20 11:15:28"
$ws1.Range("H6").Value = "2026-02-20 11:15:28"

# ---------------------------------------------------------------------
# Sheet "Estudi_Capçaleres" - source URL also points to the new period
# ---------------------------------------------------------------------
$ws2.Range("F2").Value = "https://www.meteo.cat/observacions/xema/dades?codi=XJ&dia=2026-02-20T10:30Z"
